$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget")

# The "Year 4" and "Year 5" columns (K:N) are being removed from the budget
# (the CRP now only covers Year 1-3). Before deleting the columns, rewrite
# the subtotal formulas in B/C that reference the K/L/M/N cells so that the
# column delete doesn't leave dangling #REF! errors behind - the same
# effect as if those references had never existed.
$rows = @(7, 8, 9, 11, 12, 13, 14, 15, 16)
foreach ($r in $rows) {
    $ws.Range("B$r").Formula = "=SUM(E$r,G$r,I$r)"
    $ws.Range("C$r").Formula = "=SUM(F$r,H$r,J$r)"
}

# Now remove the Year 4 / Year 5 columns outright (K:N), shifting everything
# to their right - if anything - back to the left.
$ws.Columns("K:N").Delete()

# Update the saved selection to match the post-edit state.
$ws.Range("I21").Select()
